$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 312 (D312, F312) ---
$ws.Range("D312").Value = 3.81682
$ws.Range("F312").Value = 3.7966

# --- Add new row 313 ---
$ws.Range("A312").Copy()
$ws.Range("A313").PasteSpecial(-4122)
$ws.Range("A313").Value = 45170.33333333334
$ws.Range("B313").Value = "FX_IDC:USDILS"
$ws.Range("C313").Value = 3.7966
$ws.Range("D313").Value = 3.85766
$ws.Range("E313").Value = 3.759
$ws.Range("F313").Value = 3.80432
$ws.Range("G313").Value = 0

# --- Add new row 314 ---
$ws.Range("A312").Copy()
$ws.Range("A314").PasteSpecial(-4122)
$ws.Range("A314").Value = 45201.375
$ws.Range("B314").Value = "FX_IDC:USDILS"
$ws.Range("C314").Value = 3.8155
$ws.Range("D314").Value = 4.08559
$ws.Range("E314").Value = 3.80908
$ws.Range("F314").Value = 4.0449
$ws.Range("G314").Value = 0

# --- Add new row 315 ---
$ws.Range("A312").Copy()
$ws.Range("A315").PasteSpecial(-4122)
$ws.Range("A315").Value = 45231.375
$ws.Range("B315").Value = "FX_IDC:USDILS"
$ws.Range("C315").Value = 4.0449
$ws.Range("D315").Value = 4.0449
$ws.Range("E315").Value = 3.8157
$ws.Range("F315").Value = 3.8571
$ws.Range("G315").Value = 0

$excel.CutCopyMode = 0
